# Add a new "presentations" entry for the Carpentries Instructors
# Development Meeting talk about the Data-driven CV project, and correct
# a few dates / resume-inclusion flags on the "experience" sheet that came
# up during discussion.

$wb = $excel.ActiveWorkbook

# --- presentations: new row for the "Data-driven CV" talk -----------------
$wsPres = $wb.Worksheets.Item("presentations")
$wsPres.Activate()
$wsPres.Cells.Item(3, 2).Value = "Carpentries Instructors Development Meeting"
$wsPres.Cells.Item(3, 3).Value = "December"
$wsPres.Cells.Item(3, 1).Value = "Data-driven CV"
$wsPres.Cells.Item(3, 4).Value = 2022
$wsPres.Cells.Item(3, 5).Value = "NA"
$wsPres.Cells.Item(3, 6).Value = "NA"
$wsPres.Cells.Item(3, 7).Value = "UW-Madison, WI"
$wsPres.Cells.Item(3, 8).Value = "NA"
[void]$wsPres.Range("E25").Select()

# --- experience: fix historical years and include_in_resume flags --------
$wsExp = $wb.Worksheets.Item("experience")
$wsExp.Activate()
$wsExp.Cells.Item(3, 4).Value = 2016
$wsExp.Cells.Item(5, 4).Value = 2017
$wsExp.Cells.Item(5, 9).Value = $false
$wsExp.Cells.Item(6, 4).Value = 2020
$wsExp.Cells.Item(7, 9).Value = $false
[void]$wsExp.Range("L13").Select()
